$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values cell by cell (matches upstream diff).
# For values that look like plain numbers (e.g. "1.00", "0.605"), a leading
# apostrophe forces Excel to keep them as text (preserving trailing zeros and
# exact decimal digits instead of being parsed into a float). The quote-prefix
# style flag this adds is cleared right after via Style = "Normal" so no stray
# cell-formatting changes are introduced beyond the text content itself.

$ws.Range("D2").Value = "67.016.62"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "3.266.37"
$ws.Range("E3").Value = "  +2.74%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'579.59"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").Value = "'177.23"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.267.50"
$ws.Range("E9").Value = "  +2.92%  "
$ws.Range("E10").Value = "  +4.18%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("E12").Value = "  +3.94%  "
$ws.Range("D13").Value = "3.830.43"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'28.10"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "67.011.63"
$ws.Range("E16").Value = "  +4.39%  "
$ws.Range("D17").Value = "'0.0000167"
$ws.Range("E17").Value = "  +3.14%  "
$ws.Range("D18").Value = "3.265.49"
$ws.Range("E18").Value = "  +2.69%  "
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "'13.46"
$ws.Range("D21").Value = "'371.31"
$ws.Range("E21").Value = "  +5.10%  "
$ws.Range("D22").Value = "'7.64"
$ws.Range("E22").Value = "  +5.82%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").Value = "'71.26"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").Value = "'0.512"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "3.406.46"
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  +4.68%  "
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").Value = "'22.63"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("D36").Value = "'6.81"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").Value = "'168.62"
$ws.Range("E37").Value = "  +8.64%  "
$ws.Range("E38").Value = "  +4.68%  "
$ws.Range("E39").Value = "  +6.37%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'27.28"
$ws.Range("E40").Value = "  +5.16%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.85"
$ws.Range("E41").Value = "  +9.11%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.58"
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.768.91"
$ws.Range("E43").Value = "  +4.37%  "
$ws.Range("D44").Value = "'6.44"
$ws.Range("E44").Value = "  +7.50%  "
$ws.Range("D45").Value = "'4.38"
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("D46").Value = "'346.61"
$ws.Range("E46").Value = "  +4.54%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0677"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'40.47"
$ws.Range("E48").Value = "  +4.99%  "
$ws.Range("E49").Value = "  +5.01%  "
$ws.Range("D50").Value = "'0.0279"
$ws.Range("E50").Value = "  +2.64%  "
$ws.Range("E51").Value = "  +3.32%  "

# Reset style on cells forced to text via quote-prefix so no stray formatting is introduced
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
